# Natmi following Dr Hou advice
# Rebuild rows 2-13 of Sheet1 with refreshed NATMI edge-weight stats that now
# also include the "ECs" sending cluster (previously absent from the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row: Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster, then the 16 numeric metric columns (E:T).
$rows = @(
    @("ECs", "Il4", "Il13ra2", "ECs", 1, 0.5, 0.0860255, 0.172051, 0.01866214241451323, 0.01457003365018844, 2, 1, 0.7352585, 1.470517, 0.1431143522522708, 0.1001890659257553, 0.06325098009175001, 0.253003920367, 0.002670820423292689, 0.001459758061919203),
    @("ECs", "Il4", "Il13ra2", "FAPs", 1, 0.5, 0.0860255, 0.172051, 0.01866214241451323, 0.01457003365018844, 3, 1, 4.402301, 13.206903, 0.8568856477477294, 0.8998109340742446, 0.3787101446755001, 2.272260868053, 0.01599132199122054, 0.01311027558826923),
    @("FAPs", "Il4", "Il13ra2", "ECs", 3, 1, 0.377627, 1.132881, 0.08192139369797775, 0.09593733423031037, 2, 1, 0.7352585, 1.470517, 0.1431143522522708, 0.1001890659257553, 0.2776534615795, 1.665920769477, 0.01172412719468934, 0.00961187190394179),
    @("FAPs", "Il4", "Il13ra2", "FAPs", 3, 1, 0.377627, 1.132881, 0.08192139369797775, 0.09593733423031037, 3, 1, 4.402301, 13.206903, 0.8568856477477294, 0.8998109340742446, 1.662427719727, 14.961849477543, 0.07019726650328842, 0.08632546232636858),
    @("M1", "Il4", "Il13ra2", "ECs", 3, 1, 1.258563666666667, 3.775691, 0.273029443421605, 0.3197420818403476, 2, 1, 0.7352585, 1.470517, 0.1431143522522708, 0.1001890659257553, 0.9253696337078333, 5.552217802247, 0.03907443194108101, 0.03203466051674084),
    @("M1", "Il4", "Il13ra2", "FAPs", 3, 1, 1.258563666666667, 3.775691, 0.273029443421605, 0.3197420818403476, 3, 1, 4.402301, 13.206903, 0.8568856477477294, 0.8998109340742446, 5.540576088330334, 49.865184794973, 0.233955011480524, 0.2877074213236067),
    @("M2", "Il4", "Il13ra2", "ECs", 3, 1, 0.3753126666666667, 1.125938, 0.08141932840034714, 0.09534937052400667, 2, 1, 0.7352585, 1.470517, 0.1431143522522708, 0.1001890659257553, 0.2759518283243333, 1.655710969946, 0.01165227444483059, 0.009552964369408978),
    @("M2", "Il4", "Il13ra2", "FAPs", 3, 1, 0.3753126666666667, 1.125938, 0.08141932840034714, 0.09534937052400667, 3, 1, 4.402301, 13.206903, 0.8568856477477294, 0.8998109340742446, 1.652239327779334, 14.870153950014, 0.06976705395551656, 0.08579640615459769),
    @("Neutro", "Il4", "Il13ra2", "ECs", 3, 1, 0.5777963333333334, 1.733389, 0.1253455947277287, 0.1467909867357149, 2, 1, 0.7352585, 1.470517, 0.1431143522522708, 0.1001890659257553, 0.4248296653521667, 2.548977992113, 0.01793875359713453, 0.01470685184737122),
    @("Neutro", "Il4", "Il13ra2", "FAPs", 3, 1, 0.5777963333333334, 1.733389, 0.1253455947277287, 0.1467909867357149, 3, 1, 4.402301, 13.206903, 0.8568856477477294, 0.8998109340742446, 2.543633376029667, 22.892700384267, 0.1074068411305941, 0.1320841348883437),
    @("sCs", "Il4", "Il13ra2", "ECs", 2, 1, 1.934301, 3.868602, 0.4196220973378282, 0.327610193019432, 2, 1, 0.7352585, 1.470517, 0.1431143522522708, 0.1001890659257553, 1.4222112518085, 5.688845007234001, 0.0600539446512426, 0.03282295922637331),
    @("sCs", "Il4", "Il13ra2", "FAPs", 2, 1, 1.934301, 3.868602, 0.4196220973378282, 0.327610193019432, 3, 1, 4.402301, 13.206903, 0.8568856477477294, 0.8998109340742446, 8.515375226601002, 51.092251359606, 0.3595681526865857, 0.2947872337930587)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($r, 1 + $j).Value = $rowData[$j]
    }
}

